{"js": "// The document currently has no word/styles.xml part. The target edit\n// adds that part, containing just the built-in default \"Normal\"\n// paragraph style (as happens when a styles part is (re)generated for\n// a document that only relied on the implicit default style).\n//\n// Word's object model mints/attaches the styles part the first time a\n// style is added/touched through the Styles collection, so adding the\n// (already implicitly-used) \"Normal\" paragraph style is what creates\n// word/styles.xml in the saved package.\ncontext.document.addStyle(\"Normal\", Word.StyleType.paragraph);\nawait context.sync();\n", "ps1": "# The document currently has no word/styles.xml part. The target edit\n# adds that part, containing just the built-in default \"Normal\"\n# paragraph style.\n#\n# Adding the (already implicitly-used) \"Normal\" paragraph style through\n# the Styles collection is what mints/attaches word/styles.xml in the\n# saved package.\n$d = $word.ActiveDocument\n$d.Styles.Add(\"Normal\", 1) | Out-Null\n"}
